# SYDATA.xlsx update
#  - Updated SO Header creation based on Order Number
#  - Created Test Plans
#  - Added TC for Credit Hold
#  - SO fulfillment through API
#
# This script reproduces the worksheet-level data/selection changes captured
# in the target OOXML diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AddLine: replace the 5 "a7S5f000000kRv*" test-data ids used in column G
# with the new "a7S5f000000kWR*/kWS*" ids (new Credit Hold / SO fulfillment
# test records). Re-saving these new shared-string values naturally drops the
# old unused strings and appends the new ones at the end of the shared
# string table, which is exactly what the target workbook shows.
# ---------------------------------------------------------------------------
$wsAddLine = $wb.Worksheets.Item("AddLine")
$wsAddLine.Range("G2").Value = "a7S5f000000kWRq"
$wsAddLine.Range("G3").Value = "a7S5f000000kWRv"
$wsAddLine.Range("G4").Value = "a7S5f000000kWS0"
$wsAddLine.Range("G5").Value = "a7S5f000000kWS5"
$wsAddLine.Range("G6").Value = "a7S5f000000kWSA"

# Page setup (portrait) now explicitly set on this sheet.
$wsAddLine.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# InvoiceShipper: add a "Division" column (D) used for the new Order-Number
# based SO Header creation flow.
# ---------------------------------------------------------------------------
$wsInvoiceShipper = $wb.Worksheets.Item("InvoiceShipper")
$wsInvoiceShipper.Range("D1").Value = "Division"
$wsInvoiceShipper.Range("D2").Value = "a8d5f0000004JbG"

# ---------------------------------------------------------------------------
# Batch-processing sheets: the "ProcessAllLines/Batch" flag is switched off
# (was defaulted to TRUE, now defaults to FALSE) now that individual Test
# Plans drive these flows.
# ---------------------------------------------------------------------------
$wsApproveInvoiceBatch = $wb.Worksheets.Item("ApproveInvoiceBatch")
$wsApproveInvoiceBatch.Range("B2").Value = $false

$wsDeapproveInvoiceBatch = $wb.Worksheets.Item("DeapproveInvoiceBatch")
$wsDeapproveInvoiceBatch.Range("B2").Value = $false

$wsReopenInvoiceBatch = $wb.Worksheets.Item("ReopenInvoiceBatch")
$wsReopenInvoiceBatch.Range("B2").Value = $false

$wsCloseInvoiceBatch = $wb.Worksheets.Item("CloseInvoiceBatch")
$wsCloseInvoiceBatch.Range("B2").Value = $false

$wsTransferInvoiceBatch = $wb.Worksheets.Item("TransferInvoiceBatch")
$wsTransferInvoiceBatch.Range("B2").Value = $false

# ---------------------------------------------------------------------------
# Restore the last-used selection (active cell) on every sheet that was
# touched while editing/testing the new Test Plans, then re-activate
# "AddLine" as the sheet that is on screen when the file is saved.
# ---------------------------------------------------------------------------
$wsAddHeader = $wb.Worksheets.Item("AddHeader")
$wsAddHeader.Activate()
$wsAddHeader.Range("D14").Select()

$wsInvoiceShipper.Activate()
$wsInvoiceShipper.Range("D8").Select()

$wsApproveInvoiceBatch.Activate()
$wsApproveInvoiceBatch.Range("C2").Select()

$wsDeapproveInvoiceBatch.Activate()
$wsDeapproveInvoiceBatch.Range("P14").Select()

$wsReopenInvoiceBatch.Activate()
$wsReopenInvoiceBatch.Range("H8").Select()

$wsCloseInvoiceBatch.Activate()
$wsCloseInvoiceBatch.Range("G8").Select()

$wsTransferInvoiceBatch.Activate()
$wsTransferInvoiceBatch.Range("U8").Select()

$wsPrepaymentReleaseRecall = $wb.Worksheets.Item("PrepaymentReleaseRecall")
$wsPrepaymentReleaseRecall.Activate()
$wsPrepaymentReleaseRecall.Range("F17").Select()

# AddLine becomes the active/visible sheet & selection when the workbook is
# saved.
$wsAddLine.Activate()
$wsAddLine.Range("G10").Select()
